$d = $word.ActiveDocument

# --- First paragraph: placeholder id text, trailing space cleanup, ---
# --- paragraph border spacing, and left indent.                    ---
$p1 = $d.Paragraphs.Item(1)

# 1) Swap the merge-field placeholder for the new topic id.
$p1.Range.Find.Execute("**ID__AFFARS_pgi_5315_topic_42__ID**", $false, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_SMC_PGI_5315_3G__ID**", 2)

# 2) Drop the now-orphaned trailing-space run that followed the placeholder.
$p1 = $d.Paragraphs.Item(1)
$pEnd = $p1.Range.End
$trailingSpace = $d.Range($pEnd - 2, $pEnd - 1)
if ($trailingSpace.Text -eq " ") {
    $trailingSpace.Delete()
}

# 3) Give the paragraph the same (invisible) border spacing used elsewhere
#    in the document, and widen its left indent to match.
$p1 = $d.Paragraphs.Item(1)
$p1.Borders.DistanceFromTop = 5
$p1.Borders.DistanceFromLeft = 5
$p1.Borders.DistanceFromBottom = 5
$p1.Borders.DistanceFromRight = 5
$p1.LeftIndent = 11.25

Write-Output "ok"
